# Marksheet A level - GENERAL PAPER
# Fill in "Paper 1" (column E) scores for the "Senior Five" sheet, and add
# two students (with their scores) that were missing from the sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Senior Five")

# Newly-recorded Paper 1 marks for existing students (rows 2-17).
$rowsWithMarks = @(2, 3, 4, 5, 7, 9, 13, 15, 16, 17)
$marks = @(27, 33, 45, 53, 43, 34, 30, 52, 29, 29)

for ($i = 0; $i -lt $rowsWithMarks.Length; $i++) {
    $r = $rowsWithMarks[$i]
    $cell = $ws1.Cells.Item($r, 5)
    $cell.Value = $marks[$i]
    # Nudge the cell formatting so it forks from the shared "blank" style,
    # matching how the sheet now distinguishes filled-in mark cells.
    $cell.WrapText = $false
}

# Two students who were missing from the roster entirely - add them with
# their Paper 1 marks.
$newStudents = @("OCHORA IVAN", "PILOYA MERCY")
$newMarks = @(53, 33)
$newRows = @(18, 19)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $newRows[$i]

    $nameCell = $ws1.Cells.Item($r, 2)
    $nameCell.Value = $newStudents[$i]
    $nameCell.Font.ThemeColor = 1
    $nameCell.WrapText = $false

    $markCell = $ws1.Cells.Item($r, 5)
    $markCell.Value = $newMarks[$i]
    $markCell.Font.ThemeColor = 1
    $markCell.WrapText = $false
}
